$wb = $excel.ActiveWorkbook

# --- Transformer sheet: update regulator tap positions ---
$tx = $wb.Worksheets.Item("Transformer")
$tx.Range("P21").Value = 3
$tx.Range("Q21").Value = 3
$tx.Range("R21").Value = 3
$tx.Range("P23").Value = -1
$tx.Range("P24").Value = -2
$tx.Range("P26").Value = 4
$tx.Range("P28").Value = 8

# --- Bus sheet: re-order phase rows within each 3-phase (or 2-phase) bus group ---
# Each group of consecutive C,A,B (or A,B) rows is cyclically shifted down by one
# row (last row of the group wraps around to become the first row of the group),
# i.e. new_row[i] = old_row[i-1] for i in the group, with wraparound.
$bus = $wb.Worksheets.Item("Bus")

$bus.Range("A2").Value = "1_B"
$bus.Range("E2").Value = -120
$bus.Range("A3").Value = "1_C"
$bus.Range("E3").Value = 120
$bus.Range("A4").Value = "1_A"
$bus.Range("E4").Value = 0
$bus.Range("A6").Value = "100_B"
$bus.Range("E6").Value = -120
$bus.Range("A7").Value = "100_C"
$bus.Range("E7").Value = 120
$bus.Range("A8").Value = "100_A"
$bus.Range("E8").Value = 0
$bus.Range("A9").Value = "101_B"
$bus.Range("E9").Value = -120
$bus.Range("A10").Value = "101_C"
$bus.Range("E10").Value = 120
$bus.Range("A11").Value = "101_A"
$bus.Range("E11").Value = 0
$bus.Range("A15").Value = "105_B"
$bus.Range("E15").Value = -120
$bus.Range("A16").Value = "105_C"
$bus.Range("E16").Value = 120
$bus.Range("A17").Value = "105_A"
$bus.Range("E17").Value = 0
$bus.Range("A20").Value = "108_B"
$bus.Range("E20").Value = -120
$bus.Range("A21").Value = "108_C"
$bus.Range("E21").Value = 120
$bus.Range("A22").Value = "108_A"
$bus.Range("E22").Value = 0
$bus.Range("A31").Value = "13_B"
$bus.Range("E31").Value = -120
$bus.Range("A32").Value = "13_C"
$bus.Range("E32").Value = 120
$bus.Range("A33").Value = "13_A"
$bus.Range("E33").Value = 0
$bus.Range("A34").Value = "135_B"
$bus.Range("E34").Value = -120
$bus.Range("A35").Value = "135_C"
$bus.Range("E35").Value = 120
$bus.Range("A36").Value = "135_A"
$bus.Range("E36").Value = 0
$bus.Range("A38").Value = "149_B"
$bus.Range("E38").Value = -120
$bus.Range("A39").Value = "149_C"
$bus.Range("E39").Value = 120
$bus.Range("A40").Value = "149_A"
$bus.Range("E40").Value = 0
$bus.Range("A42").Value = "150_B"
$bus.Range("E42").Value = -120
$bus.Range("A43").Value = "150_C"
$bus.Range("E43").Value = 120
$bus.Range("A44").Value = "150_A"
$bus.Range("E44").Value = 0
$bus.Range("A45").Value = "150r_B"
$bus.Range("E45").Value = -120
$bus.Range("A46").Value = "150r_C"
$bus.Range("E46").Value = 120
$bus.Range("A47").Value = "150r_A"
$bus.Range("E47").Value = 0
$bus.Range("A48").Value = "151_B"
$bus.Range("E48").Value = -120
$bus.Range("A49").Value = "151_C"
$bus.Range("E49").Value = 120
$bus.Range("A50").Value = "151_A"
$bus.Range("E50").Value = 0
$bus.Range("A51").Value = "152_B"
$bus.Range("E51").Value = -120
$bus.Range("A52").Value = "152_C"
$bus.Range("E52").Value = 120
$bus.Range("A53").Value = "152_A"
$bus.Range("E53").Value = 0
$bus.Range("A55").Value = "160_B"
$bus.Range("E55").Value = -120
$bus.Range("A56").Value = "160_C"
$bus.Range("E56").Value = 120
$bus.Range("A57").Value = "160_A"
$bus.Range("E57").Value = 0
$bus.Range("A58").Value = "160r_B"
$bus.Range("E58").Value = -120
$bus.Range("A59").Value = "160r_C"
$bus.Range("E59").Value = 120
$bus.Range("A60").Value = "160r_A"
$bus.Range("E60").Value = 0
$bus.Range("A62").Value = "18_B"
$bus.Range("E62").Value = -120
$bus.Range("A63").Value = "18_C"
$bus.Range("E63").Value = 120
$bus.Range("A64").Value = "18_A"
$bus.Range("E64").Value = 0
$bus.Range("A66").Value = "197_B"
$bus.Range("E66").Value = -120
$bus.Range("A67").Value = "197_C"
$bus.Range("E67").Value = 120
$bus.Range("A68").Value = "197_A"
$bus.Range("E68").Value = 0
$bus.Range("A71").Value = "21_B"
$bus.Range("E71").Value = -120
$bus.Range("A72").Value = "21_C"
$bus.Range("E72").Value = 120
$bus.Range("A73").Value = "21_A"
$bus.Range("E73").Value = 0
$bus.Range("A75").Value = "23_B"
$bus.Range("E75").Value = -120
$bus.Range("A76").Value = "23_C"
$bus.Range("E76").Value = 120
$bus.Range("A77").Value = "23_A"
$bus.Range("E77").Value = 0
$bus.Range("A79").Value = "25_B"
$bus.Range("E79").Value = -120
$bus.Range("A80").Value = "25_C"
$bus.Range("E80").Value = 120
$bus.Range("A81").Value = "25_A"
$bus.Range("E81").Value = 0
$bus.Range("A82").Value = "250_B"
$bus.Range("E82").Value = -120
$bus.Range("A83").Value = "250_C"
$bus.Range("E83").Value = 120
$bus.Range("A84").Value = "250_A"
$bus.Range("E84").Value = 0
$bus.Range("A91").Value = "28_B"
$bus.Range("E91").Value = -120
$bus.Range("A92").Value = "28_C"
$bus.Range("E92").Value = 120
$bus.Range("A93").Value = "28_A"
$bus.Range("E93").Value = 0
$bus.Range("A94").Value = "29_B"
$bus.Range("E94").Value = -120
$bus.Range("A95").Value = "29_C"
$bus.Range("E95").Value = 120
$bus.Range("A96").Value = "29_A"
$bus.Range("E96").Value = 0
$bus.Range("A98").Value = "30_B"
$bus.Range("E98").Value = -120
$bus.Range("A99").Value = "30_C"
$bus.Range("E99").Value = 120
$bus.Range("A100").Value = "30_A"
$bus.Range("E100").Value = 0
$bus.Range("A101").Value = "300_B"
$bus.Range("E101").Value = -120
$bus.Range("A102").Value = "300_C"
$bus.Range("E102").Value = 120
$bus.Range("A103").Value = "300_A"
$bus.Range("E103").Value = 0
$bus.Range("A104").Value = "300_open_B"
$bus.Range("E104").Value = -120
$bus.Range("A105").Value = "300_open_C"
$bus.Range("E105").Value = 120
$bus.Range("A106").Value = "300_open_A"
$bus.Range("E106").Value = 0
$bus.Range("A111").Value = "35_B"
$bus.Range("E111").Value = -120
$bus.Range("A112").Value = "35_C"
$bus.Range("E112").Value = 120
$bus.Range("A113").Value = "35_A"
$bus.Range("E113").Value = 0
$bus.Range("A114").Value = "36_B"
$bus.Range("E114").Value = -120
$bus.Range("A115").Value = "36_A"
$bus.Range("E115").Value = 0
$bus.Range("A120").Value = "40_B"
$bus.Range("E120").Value = -120
$bus.Range("A121").Value = "40_C"
$bus.Range("E121").Value = 120
$bus.Range("A122").Value = "40_A"
$bus.Range("E122").Value = 0
$bus.Range("A124").Value = "42_B"
$bus.Range("E124").Value = -120
$bus.Range("A125").Value = "42_C"
$bus.Range("E125").Value = 120
$bus.Range("A126").Value = "42_A"
$bus.Range("E126").Value = 0
$bus.Range("A128").Value = "44_B"
$bus.Range("E128").Value = -120
$bus.Range("A129").Value = "44_C"
$bus.Range("E129").Value = 120
$bus.Range("A130").Value = "44_A"
$bus.Range("E130").Value = 0
$bus.Range("A132").Value = "450_B"
$bus.Range("E132").Value = -120
$bus.Range("A133").Value = "450_C"
$bus.Range("E133").Value = 120
$bus.Range("A134").Value = "450_A"
$bus.Range("E134").Value = 0
$bus.Range("A136").Value = "47_B"
$bus.Range("E136").Value = -120
$bus.Range("A137").Value = "47_C"
$bus.Range("E137").Value = 120
$bus.Range("A138").Value = "47_A"
$bus.Range("E138").Value = 0
$bus.Range("A139").Value = "48_B"
$bus.Range("E139").Value = -120
$bus.Range("A140").Value = "48_C"
$bus.Range("E140").Value = 120
$bus.Range("A141").Value = "48_A"
$bus.Range("E141").Value = 0
$bus.Range("A142").Value = "49_B"
$bus.Range("E142").Value = -120
$bus.Range("A143").Value = "49_C"
$bus.Range("E143").Value = 120
$bus.Range("A144").Value = "49_A"
$bus.Range("E144").Value = 0
$bus.Range("A146").Value = "50_B"
$bus.Range("E146").Value = -120
$bus.Range("A147").Value = "50_C"
$bus.Range("E147").Value = 120
$bus.Range("A148").Value = "50_A"
$bus.Range("E148").Value = 0
$bus.Range("A149").Value = "51_B"
$bus.Range("E149").Value = -120
$bus.Range("A150").Value = "51_C"
$bus.Range("E150").Value = 120
$bus.Range("A151").Value = "51_A"
$bus.Range("E151").Value = 0
$bus.Range("A152").Value = "52_B"
$bus.Range("E152").Value = -120
$bus.Range("A153").Value = "52_C"
$bus.Range("E153").Value = 120
$bus.Range("A154").Value = "52_A"
$bus.Range("E154").Value = 0
$bus.Range("A155").Value = "53_B"
$bus.Range("E155").Value = -120
$bus.Range("A156").Value = "53_C"
$bus.Range("E156").Value = 120
$bus.Range("A157").Value = "53_A"
$bus.Range("E157").Value = 0
$bus.Range("A158").Value = "54_B"
$bus.Range("E158").Value = -120
$bus.Range("A159").Value = "54_C"
$bus.Range("E159").Value = 120
$bus.Range("A160").Value = "54_A"
$bus.Range("E160").Value = 0
$bus.Range("A161").Value = "55_B"
$bus.Range("E161").Value = -120
$bus.Range("A162").Value = "55_C"
$bus.Range("E162").Value = 120
$bus.Range("A163").Value = "55_A"
$bus.Range("E163").Value = 0
$bus.Range("A164").Value = "56_B"
$bus.Range("E164").Value = -120
$bus.Range("A165").Value = "56_C"
$bus.Range("E165").Value = 120
$bus.Range("A166").Value = "56_A"
$bus.Range("E166").Value = 0
$bus.Range("A167").Value = "57_B"
$bus.Range("E167").Value = -120
$bus.Range("A168").Value = "57_C"
$bus.Range("E168").Value = 120
$bus.Range("A169").Value = "57_A"
$bus.Range("E169").Value = 0
$bus.Range("A173").Value = "60_B"
$bus.Range("E173").Value = -120
$bus.Range("A174").Value = "60_C"
$bus.Range("E174").Value = 120
$bus.Range("A175").Value = "60_A"
$bus.Range("E175").Value = 0
$bus.Range("A176").Value = "61_B"
$bus.Range("E176").Value = -120
$bus.Range("A177").Value = "61_C"
$bus.Range("E177").Value = 120
$bus.Range("A178").Value = "61_A"
$bus.Range("E178").Value = 0
$bus.Range("A179").Value = "610_B"
$bus.Range("E179").Value = -120
$bus.Range("A180").Value = "610_C"
$bus.Range("E180").Value = 120
$bus.Range("A181").Value = "610_A"
$bus.Range("E181").Value = 0
$bus.Range("A182").Value = "61s_B"
$bus.Range("E182").Value = -120
$bus.Range("A183").Value = "61s_C"
$bus.Range("E183").Value = 120
$bus.Range("A184").Value = "61s_A"
$bus.Range("E184").Value = 0
$bus.Range("A185").Value = "62_B"
$bus.Range("E185").Value = -120
$bus.Range("A186").Value = "62_C"
$bus.Range("E186").Value = 120
$bus.Range("A187").Value = "62_A"
$bus.Range("E187").Value = 0
$bus.Range("A188").Value = "63_B"
$bus.Range("E188").Value = -120
$bus.Range("A189").Value = "63_C"
$bus.Range("E189").Value = 120
$bus.Range("A190").Value = "63_A"
$bus.Range("E190").Value = 0
$bus.Range("A191").Value = "64_B"
$bus.Range("E191").Value = -120
$bus.Range("A192").Value = "64_C"
$bus.Range("E192").Value = 120
$bus.Range("A193").Value = "64_A"
$bus.Range("E193").Value = 0
$bus.Range("A194").Value = "65_B"
$bus.Range("E194").Value = -120
$bus.Range("A195").Value = "65_C"
$bus.Range("E195").Value = 120
$bus.Range("A196").Value = "65_A"
$bus.Range("E196").Value = 0
$bus.Range("A197").Value = "66_B"
$bus.Range("E197").Value = -120
$bus.Range("A198").Value = "66_C"
$bus.Range("E198").Value = 120
$bus.Range("A199").Value = "66_A"
$bus.Range("E199").Value = 0
$bus.Range("A200").Value = "67_B"
$bus.Range("E200").Value = -120
$bus.Range("A201").Value = "67_C"
$bus.Range("E201").Value = 120
$bus.Range("A202").Value = "67_A"
$bus.Range("E202").Value = 0
$bus.Range("A205").Value = "7_B"
$bus.Range("E205").Value = -120
$bus.Range("A206").Value = "7_C"
$bus.Range("E206").Value = 120
$bus.Range("A207").Value = "7_A"
$bus.Range("E207").Value = 0
$bus.Range("A210").Value = "72_B"
$bus.Range("E210").Value = -120
$bus.Range("A211").Value = "72_C"
$bus.Range("E211").Value = 120
$bus.Range("A212").Value = "72_A"
$bus.Range("E212").Value = 0
$bus.Range("A216").Value = "76_B"
$bus.Range("E216").Value = -120
$bus.Range("A217").Value = "76_C"
$bus.Range("E217").Value = 120
$bus.Range("A218").Value = "76_A"
$bus.Range("E218").Value = 0
$bus.Range("A219").Value = "77_B"
$bus.Range("E219").Value = -120
$bus.Range("A220").Value = "77_C"
$bus.Range("E220").Value = 120
$bus.Range("A221").Value = "77_A"
$bus.Range("E221").Value = 0
$bus.Range("A222").Value = "78_B"
$bus.Range("E222").Value = -120
$bus.Range("A223").Value = "78_C"
$bus.Range("E223").Value = 120
$bus.Range("A224").Value = "78_A"
$bus.Range("E224").Value = 0
$bus.Range("A225").Value = "79_B"
$bus.Range("E225").Value = -120
$bus.Range("A226").Value = "79_C"
$bus.Range("E226").Value = 120
$bus.Range("A227").Value = "79_A"
$bus.Range("E227").Value = 0
$bus.Range("A228").Value = "8_B"
$bus.Range("E228").Value = -120
$bus.Range("A229").Value = "8_C"
$bus.Range("E229").Value = 120
$bus.Range("A230").Value = "8_A"
$bus.Range("E230").Value = 0
$bus.Range("A231").Value = "80_B"
$bus.Range("E231").Value = -120
$bus.Range("A232").Value = "80_C"
$bus.Range("E232").Value = 120
$bus.Range("A233").Value = "80_A"
$bus.Range("E233").Value = 0
$bus.Range("A234").Value = "81_B"
$bus.Range("E234").Value = -120
$bus.Range("A235").Value = "81_C"
$bus.Range("E235").Value = 120
$bus.Range("A236").Value = "81_A"
$bus.Range("E236").Value = 0
$bus.Range("A237").Value = "82_B"
$bus.Range("E237").Value = -120
$bus.Range("A238").Value = "82_C"
$bus.Range("E238").Value = 120
$bus.Range("A239").Value = "82_A"
$bus.Range("E239").Value = 0
$bus.Range("A240").Value = "83_B"
$bus.Range("E240").Value = -120
$bus.Range("A241").Value = "83_C"
$bus.Range("E241").Value = 120
$bus.Range("A242").Value = "83_A"
$bus.Range("E242").Value = 0
$bus.Range("A245").Value = "86_B"
$bus.Range("E245").Value = -120
$bus.Range("A246").Value = "86_C"
$bus.Range("E246").Value = 120
$bus.Range("A247").Value = "86_A"
$bus.Range("E247").Value = 0
$bus.Range("A248").Value = "87_B"
$bus.Range("E248").Value = -120
$bus.Range("A249").Value = "87_C"
$bus.Range("E249").Value = 120
$bus.Range("A250").Value = "87_A"
$bus.Range("E250").Value = 0
$bus.Range("A252").Value = "89_B"
$bus.Range("E252").Value = -120
$bus.Range("A253").Value = "89_C"
$bus.Range("E253").Value = 120
$bus.Range("A254").Value = "89_A"
$bus.Range("E254").Value = 0
$bus.Range("A257").Value = "91_B"
$bus.Range("E257").Value = -120
$bus.Range("A258").Value = "91_C"
$bus.Range("E258").Value = 120
$bus.Range("A259").Value = "91_A"
$bus.Range("E259").Value = 0
$bus.Range("A261").Value = "93_B"
$bus.Range("E261").Value = -120
$bus.Range("A262").Value = "93_C"
$bus.Range("E262").Value = 120
$bus.Range("A263").Value = "93_A"
$bus.Range("E263").Value = 0
$bus.Range("A266").Value = "95_B"
$bus.Range("E266").Value = -120
$bus.Range("A267").Value = "95_C"
$bus.Range("E267").Value = 120
$bus.Range("A268").Value = "95_A"
$bus.Range("E268").Value = 0
$bus.Range("A270").Value = "97_B"
$bus.Range("E270").Value = -120
$bus.Range("A271").Value = "97_C"
$bus.Range("E271").Value = 120
$bus.Range("A272").Value = "97_A"
$bus.Range("E272").Value = 0
$bus.Range("A273").Value = "98_B"
$bus.Range("E273").Value = -120
$bus.Range("A274").Value = "98_C"
$bus.Range("E274").Value = 120
$bus.Range("A275").Value = "98_A"
$bus.Range("E275").Value = 0
$bus.Range("A276").Value = "99_B"
$bus.Range("E276").Value = -120
$bus.Range("A277").Value = "99_C"
$bus.Range("E277").Value = 120
$bus.Range("A278").Value = "99_A"
$bus.Range("E278").Value = 0
